# Apply updated cryptocurrency price/volume data and coin reorderings
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.635.00'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '1.754.07'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.19'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4478'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.42%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3541'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07418'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.79'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.087'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.47%  '
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.72'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.984'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.156'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.72%  '
$ws.Range('D16').Value = '1.751.63'
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.52'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001057'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06423'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.04'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.742'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.57%  '
$ws.Range('D23').Value = '27.683.57'
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.21'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.120'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.37'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.21'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').Value = '1.955.59'
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.074'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.24%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.44'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.062'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09138'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.71%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.661'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.461'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02284'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '11.68'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.40%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06041'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2067'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.95%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6283'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.98%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.942'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.183'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.86%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.381'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.769'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.18'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.04%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5879'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.52%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.710'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '122.90'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.947'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06906'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.121'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.49'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.25%  '
